# Add Portfolioüberblick rows for M4 (rows 30-41) to Tabelle1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row data: A=Meilenstein, B=Dokument, C=Kapitel, D=Beschreibung, E..H=Anteile
$rows = @(
    @{ Row=30; A=4; B="Aktivitätsdiagramm.vsd";          D="Aktivtätsdiagramm";              E=10; F=10; G=10; H=70 },
    @{ Row=31; A=4; B="Klassendiagram.vsd";               D="Klassendiagramm";                E=10; F=70; G=10; H=10 },
    @{ Row=32; A=4; B="Komponentendiagramm.vsd";          D="Komponentendiagramm";            E=30; F=30; G=30; H=10 },
    @{ Row=33; A=4; B="seq_dia_reservation_v2.vsd";       D="Sequenzdiagramm";                E=10; F=10; G=70; H=10 },
    @{ Row=34; A=4; B="Zustandsdiagramm.vsd";             D="Zustandsdiagramm";               E=30; F=30; G=30; H=10 },
    @{ Row=35; A=4; B="FST17_M4_Aerzteapp.doc"; C="1 Architekturkonzept";                     E=70; F=10; G=10; H=10 },
    @{ Row=36; A=4; B="FST17_M4_Aerzteapp.doc"; C="2 Komponentendiagramm";                    E=10; F=70; G=10; H=10 },
    @{ Row=37; A=4; B="FST17_M4_Aerzteapp.doc"; C="3 Grobentwurf (Klassendiagramm)";           E=10; F=70; G=10; H=10 },
    @{ Row=38; A=4; B="FST17_M4_Aerzteapp.doc"; C="4 Sequenzdiagramm";                         E=10; F=10; G=70; H=10 },
    @{ Row=39; A=4; B="FST17_M4_Aerzteapp.doc"; C="5 Aktivitätsdiagramm";                      E=10; F=10; G=10; H=70 },
    @{ Row=40; A=4; B="FST17_M4_Aerzteapp.doc"; C="6 Zustandsdiagramm";                        E=70; F=10; G=10; H=10 },
    @{ Row=41; A=4; B="FST17_M4_Aerzteapp.doc"; C="7 Realisierungstechnologien";               E=70; F=10; G=10; H=10 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    if ($r.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $r.C
    }
    if ($r.ContainsKey("D")) {
        $ws.Cells.Item($row, 4).Value = $r.D
    }
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Formula = "=SUM(E" + $row + ":H" + $row + ")"
}

# Update view: scroll down and select the last filled cell, matching the
# post-edit author selection state.
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("I41").Select()
